$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are numeric-looking strings that must be
# preserved exactly as typed (e.g. trailing zeros: "1.00", "0.790",
# "0.0900") -- format as Text first so Excel does not normalize them
# into plain numbers.
$forceTextCells = @("D23", "D40", "D46")
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '56.034.26'
$ws.Range("E2").Value = '  -3.12%  '
$ws.Range("D3").Value = '2.361.95'
$ws.Range("E3").Value = '  -3.90%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '500.66'
$ws.Range("D6").Value = '129.23'
$ws.Range("E6").Value = '  -3.48%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '0.543'
$ws.Range("E8").Value = '  -2.54%  '
$ws.Range("D9").Value = '2.362.83'
$ws.Range("E9").Value = '  -3.89%  '
$ws.Range("D10").Value = '0.0981'
$ws.Range("E10").Value = '  -0.25%  '
$ws.Range("E11").Value = '  +0.18%  '
$ws.Range("D12").Value = '4.77'
$ws.Range("E12").Value = '  +2.84%  '
$ws.Range("D13").Value = '0.323'
$ws.Range("E13").Value = '  -0.32%  '
$ws.Range("D14").Value = '2.780.82'
$ws.Range("E14").Value = '  -3.82%  '
$ws.Range("D15").Value = '56.016.97'
$ws.Range("E15").Value = '  -3.28%  '
$ws.Range("D16").Value = '21.38'
$ws.Range("E16").Value = '  -2.68%  '
$ws.Range("E17").Value = '  -2.22%  '
$ws.Range("D18").Value = '2.391.20'
$ws.Range("E18").Value = '  -3.91%  '
$ws.Range("E19").Value = '  -3.60%  '
$ws.Range("E20").Value = '  -3.44%  '
$ws.Range("D21").Value = '306.56'
$ws.Range("E21").Value = '  -2.81%  '
$ws.Range("E22").Value = '  -2.93%  '
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").Value = '65.42'
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("D26").Value = '0.368'
$ws.Range("E26").Value = '  -3.37%  '
$ws.Range("E27").Value = '  -6.05%  '
$ws.Range("D28").Value = '7.21'
$ws.Range("E28").Value = '  -4.82%  '
$ws.Range("D29").Value = '171.45'
$ws.Range("E29").Value = '  -0.67%  '
$ws.Range("D30").Value = '0.0₃0710'
$ws.Range("E30").Value = '  -3.37%  '
$ws.Range("D31").Value = '1.64'
$ws.Range("E31").Value = '  -3.34%  '
$ws.Range("E32").Value = '  +0.03%  '
$ws.Range("B33").Value = 'Aptos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D33").Value = '5.74'
$ws.Range("E33").Value = '  -7.08%  '
$ws.Range("B34").Value = 'FirstDigitalUSD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D34").Value = '0.997'
$ws.Range("E34").Value = '  -0.12%  '
$ws.Range("E35").Value = '  -4.91%  '
$ws.Range("E36").Value = '  -3.04%  '
$ws.Range("E37").Value = '  -6.04%  '
$ws.Range("E38").Value = '  -3.94%  '
$ws.Range("D39").Value = '36.08'
$ws.Range("E39").Value = '  -1.89%  '
$ws.Range("D40").Value = '0.790'
$ws.Range("E41").Value = '  -5.92%  '
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").Value = '3.35'
$ws.Range("E42").Value = '  -2.20%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").Value = '128.61'
$ws.Range("E43").Value = '  -6.24%  '
$ws.Range("D44").Value = '4.71'
$ws.Range("E44").Value = '  -4.14%  '
$ws.Range("D45").Value = '0.561'
$ws.Range("E45").Value = '  -2.99%  '
$ws.Range("D46").Value = '0.0900'
$ws.Range("E46").Value = '  -2.28%  '
$ws.Range("D47").Value = '238.99'
$ws.Range("E47").Value = '  -6.93%  '
$ws.Range("E48").Value = '  -2.78%  '
$ws.Range("E49").Value = '  -3.98%  '
$ws.Range("D50").Value = '17.04'
$ws.Range("E50").Value = '  -2.21%  '
$ws.Range("E51").Value = '  -0.57%  '
